$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = [double]"8.191115064622072e-18"
$ws.Cells.Item(3, 2).Value = [double]"0"
$ws.Cells.Item(4, 2).Value = [double]"2.637861236617057e-18"
$ws.Cells.Item(5, 2).Value = [double]"0"
$ws.Cells.Item(6, 2).Value = [double]"0"
$ws.Cells.Item(7, 2).Value = [double]"0"
$ws.Cells.Item(8, 2).Value = [double]"0"
$ws.Cells.Item(9, 2).Value = [double]"0"
$ws.Cells.Item(10, 2).Value = [double]"0"
$ws.Cells.Item(11, 2).Value = [double]"0.0293397299680501"
$ws.Cells.Item(12, 2).Value = [double]"9.964516315928008e-19"
$ws.Cells.Item(13, 2).Value = [double]"2.206582755816914e-18"
$ws.Cells.Item(14, 2).Value = [double]"1.335138492205378e-17"
$ws.Cells.Item(15, 2).Value = [double]"0.1242669696252968"
$ws.Cells.Item(16, 2).Value = [double]"3.377615051119878e-18"
$ws.Cells.Item(17, 2).Value = [double]"0"
$ws.Cells.Item(18, 2).Value = [double]"0"
$ws.Cells.Item(19, 2).Value = [double]"0.07035691606082364"
$ws.Cells.Item(20, 2).Value = [double]"0"
$ws.Cells.Item(21, 2).Value = [double]"0"
$ws.Cells.Item(22, 2).Value = [double]"0"
$ws.Cells.Item(23, 2).Value = [double]"0"
$ws.Cells.Item(24, 2).Value = [double]"0.1121593627377357"
$ws.Cells.Item(25, 2).Value = [double]"2.81741002183871e-18"
$ws.Cells.Item(26, 2).Value = [double]"0"
$ws.Cells.Item(27, 2).Value = [double]"0"
$ws.Cells.Item(28, 2).Value = [double]"0.1078670493821004"
$ws.Cells.Item(29, 2).Value = [double]"4.068759754403245e-18"
$ws.Cells.Item(30, 2).Value = [double]"0"
$ws.Cells.Item(31, 2).Value = [double]"0"
$ws.Cells.Item(32, 2).Value = [double]"5.889486739625472e-18"
$ws.Cells.Item(33, 2).Value = [double]"0.02262602944344694"
$ws.Cells.Item(34, 2).Value = [double]"4.663083637972925e-18"
$ws.Cells.Item(35, 2).Value = [double]"0.03526033269644029"
$ws.Cells.Item(36, 2).Value = [double]"0.01415084953189176"
$ws.Cells.Item(37, 2).Value = [double]"0.1189297524597908"
$ws.Cells.Item(38, 2).Value = [double]"0.06596603463736081"
$ws.Cells.Item(39, 2).Value = [double]"2.221921836205665e-17"
$ws.Cells.Item(40, 2).Value = [double]"0.04933318210257538"
$ws.Cells.Item(41, 2).Value = [double]"0"
$ws.Cells.Item(42, 2).Value = [double]"0.13573278575077"
$ws.Cells.Item(43, 2).Value = [double]"0"
$ws.Cells.Item(44, 2).Value = [double]"0.06637694019961772"
$ws.Cells.Item(45, 2).Value = [double]"0"
$ws.Cells.Item(46, 2).Value = [double]"0.03521209812160288"
$ws.Cells.Item(47, 2).Value = [double]"0"
$ws.Cells.Item(48, 2).Value = [double]"0.01242196728249686"
$ws.Cells.Item(49, 2).Value = [double]"5.199317760728858e-18"
$ws.Cells.Item(50, 2).Value = [double]"0"
$ws.Cells.Item(51, 2).Value = [double]"0"
